# Apply cryptos list price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.783.87'
$ws.Range("E2").Value = '  +0.30%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.352.52'
$ws.Range("E3").Value = '  -0.17%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '554.88'
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.12'
$ws.Range("E6").Value = '  -1.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.629'
$ws.Range("E7").Value = '  +2.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.343.68'
$ws.Range("E8").Value = '  -0.19%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.175'
$ws.Range("E10").Value = '  +7.60%  '
$ws.Range("E11").Value = '  +1.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.66'
$ws.Range("E12").Value = '  -2.49%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000281'
$ws.Range("E13").Value = '  +3.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.10'
$ws.Range("E14").Value = '  +0.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.875.62'
$ws.Range("E15").Value = '  -0.66%  '
$ws.Range("E16").Value = '  +2.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.19'
$ws.Range("E17").Value = '  -0.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.350.00'
$ws.Range("E18").Value = '  -0.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '64.672.74'
$ws.Range("E19").Value = '  +0.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.78'
$ws.Range("E20").Value = '  -0.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.989'
$ws.Range("E21").Value = '  +0.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '448.34'
$ws.Range("E22").Value = '  +2.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.95'
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.07'
$ws.Range("E24").Value = '  -0.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '87.18'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.75'
$ws.Range("E26").Value = '  +2.89%  '
$ws.Range("E27").Value = '  +1.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.67'
$ws.Range("E28").Value = '  -1.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.64'
$ws.Range("E29").Value = '  -1.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.91'
$ws.Range("E30").Value = '  +3.95%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.53'
$ws.Range("E31").Value = '  -1.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '62.79'
$ws.Range("E32").Value = '  +7.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.42'
$ws.Range("E33").Value = '  -0.42%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '574.06'
$ws.Range("E34").Value = '  -0.96%  '
$ws.Range("E35").Value = '  -0.52%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("B37").Value = 'Stacks'
$ws.Range("C37").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.59'
$ws.Range("E37").Value = '  +2.11%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.142'
$ws.Range("E38").Value = '  -0.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.57'
$ws.Range("E39").Value = '  -0.46%  '
$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0745'
$ws.Range("E40").Value = '  -1.30%  '
$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.369'
$ws.Range("E41").Value = '  +0.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.085.78'
$ws.Range("E42").Value = '  -0.66%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0417'
$ws.Range("E43").Value = '  +1.40%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.74'
$ws.Range("E44").Value = '  -2.69%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.134'
$ws.Range("E45").Value = '  +3.27%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.17'
$ws.Range("E46").Value = '  -1.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.44'
$ws.Range("E47").Value = '  -0.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  +0.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '140.63'
$ws.Range("E49").Value = '  +4.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.51'
$ws.Range("E50").Value = '  -2.62%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.27'
$ws.Range("E51").Value = '  -0.43%  '
